$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column),
# shifting the existing "Late", "Date", "Loan Balance" columns
# (N, O, P) one place to the right (to O, P, Q).
$ws.Columns("N").Insert()

# The new column takes on the width of the column to its left (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selected range on the sheet to match the new state.
$ws.Range("S7").Select()
